# Apply "syntax sugar" edits to the scenario sheet of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario")

$ws.Range("G7").Value = 'confirmation_request="Sorry to ask you again, but do you often eat sadwiches?"'
$ws.Range("F11").Value = 'confirmed=="yes"'
$ws.Range("F12").Value = 'confirmed=="no"'
$ws.Range("F15").Value = '$"Please determine if the user said the reason."'
$ws.Range("F17").Value = '#favorite-sandwich=="egg salad sandwich"'
$ws.Range("G17").Value = 'topic_sandwich=#favorite-sandwich'
$ws.Range("G18").Value = 'topic_sandwich=#favorite-sandwich'
$ws.Range("C24").Value = "I understand. {`$`"Generate a sentence to say it's time to end the talk by continuing the conversation in 50 words`" }  Thank you for your time."
$ws.Range("B25").Value = "#final_apology"

# Selection / scroll position change recorded in the sheet view.
$ws.Range("F15").Select()
$excel.ActiveWindow.ScrollRow = 6
